$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 2000
$ws.Range("K74").Value = 2000
$ws.Range("M74").Value = -1064

# Row 77
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 2000
$ws.Range("K77").Value = 10000
$ws.Range("M77").Value = -5320

$ws = $wb.Worksheets.Item("ARM")
# Row 76
$ws.Range("H76").Value = 23657.4
$ws.Range("J76").Value = 23657.4
$ws.Range("L76").Value = 23657.4
$ws.Range("N76").Value = -24333.4

# Row 79
$ws.Range("H79").Value = 23657.4
$ws.Range("J79").Value = 23657.4
$ws.Range("L79").Value = 23657.4
$ws.Range("N79").Value = -25997.4

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2769.4
$ws.Range("I20").Value = 2616.6667
$ws.Range("K20").Value = 2616.6667
$ws.Range("M20").Value = -2369.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 900
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 1500
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = -126
$ws.Range("N17").Value = -1848

# Row 22
$ws.Range("H22").Value = 360.5
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 360.5
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 360.5
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1060.5

# Row 25
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Row 41
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4572

# Row 43
$ws.Range("H43").Value = 19663.334
$ws.Range("J43").Value = 19663.334
$ws.Range("L43").Value = 19663.334
$ws.Range("N43").Value = -20031.334

# Row 50
$ws.Range("H50").Value = 83
$ws.Range("I50").Value = 83
$ws.Range("K50").Value = 83
$ws.Range("M50").Value = 542

# Row 86
$ws.Range("H86").Value = 12098.8
$ws.Range("I86").Value = 11831.667
$ws.Range("J86").Value = 12499.5
$ws.Range("K86").Value = 11831.667
$ws.Range("L86").Value = 12499.5
$ws.Range("M86").Value = -10708.667
$ws.Range("N86").Value = -14745.5

# Row 89
$ws.Range("H89").Value = 12098.8
$ws.Range("I89").Value = 11831.667
$ws.Range("J89").Value = 12499.5
$ws.Range("K89").Value = 59158.335
$ws.Range("L89").Value = 62497.5
$ws.Range("M89").Value = -53542.335
$ws.Range("N89").Value = -73729.5

# Row 92
$ws.Range("H92").Value = 37600
$ws.Range("J92").Value = 37600
$ws.Range("L92").Value = 37600
$ws.Range("N92").Value = -42592

# Row 101
$ws.Range("H101").Value = 19663.334
$ws.Range("J101").Value = 19663.334
$ws.Range("L101").Value = 19663.334
$ws.Range("N101").Value = -26153.334

$ws = $wb.Worksheets.Item("CUL")
# Row 140
$ws.Range("H140").Value = 1004.2
$ws.Range("I140").Value = 1004.2
$ws.Range("K140").Value = 3012.6
$ws.Range("M140").Value = 2167.4

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 13620000
$ws.Range("I11").Value = 17000000
$ws.Range("J11").Value = 100000
$ws.Range("K11").Value = 17000000
$ws.Range("L11").Value = 100000
$ws.Range("M11").Value = -16999861
$ws.Range("N11").Value = -100278

# Row 80
$ws.Range("H80").Value = 1083.4
$ws.Range("J80").Value = 2325
$ws.Range("L80").Value = 2325
$ws.Range("N80").Value = -4321

# Row 83
$ws.Range("H83").Value = 1083.4
$ws.Range("J83").Value = 2325
$ws.Range("L83").Value = 11625
$ws.Range("N83").Value = -21609

# Row 97
$ws.Range("H97").Value = 4010
$ws.Range("I97").Value = 4010
$ws.Range("K97").Value = 4010
$ws.Range("M97").Value = -3514

# Row 99
$ws.Range("H99").Value = 8833.333000000001
$ws.Range("I99").Value = 8833.333000000001
$ws.Range("K99").Value = 8833.333000000001
$ws.Range("M99").Value = -6587.333000000001

# Row 132
$ws.Range("H132").Value = 6135.375
$ws.Range("I132").Value = 6135.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 18406.125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15876.125
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 4000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -4224

# Row 25
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# Row 46
$ws.Range("H46").Value = 5158
$ws.Range("I46").Value = 5529.3335
$ws.Range("J46").Value = 4786.6665
$ws.Range("K46").Value = 5529.3335
$ws.Range("L46").Value = 4786.6665
$ws.Range("M46").Value = -5341.3335
$ws.Range("N46").Value = -5162.6665

# Row 68
$ws.Range("H68").Value = 2537.6
$ws.Range("I68").Value = 2172
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2172
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1423
$ws.Range("N68").Value = -5498

# Row 71
$ws.Range("H71").Value = 2537.6
$ws.Range("I71").Value = 2172
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 10860
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -7116
$ws.Range("N71").Value = -27488

# Row 93
$ws.Range("H93").Value = 1333.4445
$ws.Range("I93").Value = 1450.125
$ws.Range("J93").Value = 400
$ws.Range("K93").Value = 1450.125
$ws.Range("L93").Value = 400
$ws.Range("M93").Value = -202.125
$ws.Range("N93").Value = -2896

# Row 126
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 12000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 34855.5
$ws.Range("I45").Value = 34711
$ws.Range("K45").Value = 34711
$ws.Range("M45").Value = -34220

# Row 80
$ws.Range("H80").Value = 40000
$ws.Range("J80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("N80").Value = -41996

# Row 83
$ws.Range("H83").Value = 40000
$ws.Range("J83").Value = 40000
$ws.Range("L83").Value = 120000
$ws.Range("N83").Value = -129984

# Row 105
$ws.Range("H105").Value = 39966
$ws.Range("J105").Value = 39966
$ws.Range("L105").Value = 39966
$ws.Range("N105").Value = -46954
